{"js": "// Remove the trailing \"Ver no Jupiter...\" line, the \"\u00a9 2020 ...\" footer\n// line, and the blank paragraph that separates them from the bibliography\n// entry above, matching the site rebuild diff that dropped the Jekyll\n// page-footer boilerplate from the end of the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the bibliography paragraph that stays (JEWETT/SERWAY reference) and\n// the two footer paragraphs that must go, by their literal text so the\n// script is resilient to any pre-existing index drift.\nlet jewettIndex = -1;\nlet jupiterIndex = -1;\nlet copyrightIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n  if (text.indexOf(\"JEWETT Jr, John W.\") !== -1) {\n    jewettIndex = i;\n  } else if (text.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIndex = i;\n  } else if (text.indexOf(\"Powered by Jekyll\") !== -1) {\n    copyrightIndex = i;\n  }\n}\n\n// Delete the footer paragraphs (copyright line, jupiter line, and the blank\n// paragraph right after the bibliography entry) in reverse document order so\n// earlier deletions don't shift the indices of paragraphs still to remove.\nif (copyrightIndex !== -1) {\n  items[copyrightIndex].delete();\n}\nif (jupiterIndex !== -1) {\n  items[jupiterIndex].delete();\n}\nif (jewettIndex !== -1 && jewettIndex + 1 < items.length) {\n  // The blank paragraph immediately following the bibliography entry.\n  items[jewettIndex + 1].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" line, the \"(c) 2020 ...\" footer\n# line, and the blank paragraph that separates them from the bibliography\n# entry above, matching the site rebuild diff that dropped the Jekyll\n# page-footer boilerplate from the end of the document.\n\n$d = $word.ActiveDocument\n\n$jewettIndex = -1\n$jupiterIndex = -1\n$copyrightIndex = -1\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"*JEWETT Jr, John W.*\") {\n        $jewettIndex = $i\n    } elseif ($t -like \"*Ver no Jupiter*\") {\n        $jupiterIndex = $i\n    } elseif ($t -like \"*Powered by Jekyll*\") {\n        $copyrightIndex = $i\n    }\n}\n\n# Delete from the bottom up so earlier deletions don't shift the index of\n# paragraphs still to be removed.\nif ($copyrightIndex -ne -1) {\n    $d.Paragraphs.Item($copyrightIndex).Range.Delete()\n}\nif ($jupiterIndex -ne -1) {\n    $d.Paragraphs.Item($jupiterIndex).Range.Delete()\n}\nif ($jewettIndex -ne -1) {\n    # The blank paragraph immediately following the bibliography entry.\n    $d.Paragraphs.Item($jewettIndex + 1).Range.Delete()\n}\n"}
